$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Support importing "lqsutra" info: insert two new columns (D,E) for
# "sutra__lqsutra" / "sutra__lqsutra__name", right after the (renamed)
# "sutra__name" column, shifting the former D:I ("code" .. "remark") to F:K.

$ws.Range("D1:E1").EntireColumn.Insert()

# Give the new D:E data cells (rows 2-3) the same formatting as the
# neighbouring A:B columns on those rows (center aligned data style),
# rather than the inherited "name" column look-and-feel.
$ws.Range("A2:B3").Copy()
$ws.Range("D2:E3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row.
$ws.Range("C1").Value = "sutra__name"
$ws.Range("D1").Value = "sutra__lqsutra"
$ws.Range("E1").Value = "sutra__lqsutra__name"

# Row 2 (QL0001 / 大清三藏聖教目錄 entry) - matching lqsutra record.
$ws.Range("D2").Value = "LQ001"
$ws.Range("E2").Value = "LQ大清三藏聖教目錄"

# Row 3 (QL0003 / 大般若波羅蜜多經 entry) - matching lqsutra record.
$ws.Range("D3").Value = "LQ003"
$ws.Range("E3").Value = "LQ大般若波羅蜜多經"

# Best-effort width for the two brand-new columns (content-fit).
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

$ws.Range("D6").Select() | Out-Null
